# RB-Blessing / T-Tabelle.xlsx
# "Anforderungsliste ergänzt, T-Tabelle mit Funktionen befüllt"
#
# The "Funktionen" (functions) list in column K (rows 5-15) is replaced with
# the real list of mechanical sub-functions/components, the S4 specification
# label gets a capitalization fix, and the active selection/scroll position
# is updated to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Spezifikation S4 label: fix capitalization ("Benutzbar" -> "benutzbar") ---
$ws.Range("E3").Value = "S4-Im freien benutzbar"

# --- "Funktionen" column (K5:K15): replace placeholder list with the real one ---
$ws.Range("K5").Value  = "Befestigung Trum - Band"
$ws.Range("K6").Value  = "Welle"
$ws.Range("K7").Value  = "Wellenlagerung"
$ws.Range("K8").Value  = "Abdichtung"
$ws.Range("K9").Value  = "Motor"
$ws.Range("K10").Value = "Verbindung Welle - Trommel"
$ws.Range("K11").Value = "Trommel"
$ws.Range("K12").Value = "Kettentrieb"
$ws.Range("K13").Value = "Lagerböcke"
$ws.Range("K14").Value = "Gehäuse"
$ws.Range("K15").Value = "Befestigung auf dem Untergrund"

# --- Restore the view state (scroll position / active cell) left by the author ---
$ws.Activate()
$ws.Range("A5").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("K18").Select()
